# Refresh the "Estado de Cuenta" (account statement) detail table.
# The Periodo Mora (E), Valor Mora (F) and Salario Basico (G) columns for the
# 36 detail rows (16-51) are updated: periods now run in ascending order
# (1704 .. 2003, instead of the previous descending 2003 .. 1704) and the
# Valor Mora / Salario Basico figures are refreshed to match "parte 1" of the
# new account-statement database.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$periodoMora = @(
    "1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809",
    "1810","1811","1812","1901","1902","1903","1904","1905","1906",
    "1907","1908","1909","1910","1911","1912","2001","2002","2003"
)

$valorMora = @(
    27578,27578,27578,27578,27578,27578,27578,27578,27578,
    27578,27578,27578,27578,27578,27578,27578,27578,31249,
    31249,31249,31249,31249,31249,31249,31249,31249,31249,
    31249,31249,31249,31249,31249,31249,31249,31249,31249
)

$salarioBasico = @(
    781242,781242,781242,781242,781242,781242,781242,781242,781242,
    781242,781242,781242,781242,781242,781242,781242,781242,781242,
    781242,781242,781242,781242,781242,781242,781242,781242,781242,
    781242,781242,781242,781242,781242,781242,781242,781242,781242
)

$firstRow = 16
for ($i = 0; $i -lt $periodoMora.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periodoMora[$i]
    $ws.Cells.Item($row, 6).Value = $valorMora[$i]
    $ws.Cells.Item($row, 7).Value = $salarioBasico[$i]
}
